# Generate Report for Handoff
# Adds two new file entries (44e48f6d-... and 9108f6ff-...) to the
# localization-status workbook: one new row per entry on the "Overview"
# sheet, and one new row per entry on each of the "zh-cn" and "de-de"
# detail sheets.

$wb = $excel.ActiveWorkbook

$ghSrc  = "https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/"
$ghZh   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b58d48b4b51d881d18df2f827562167da1273289/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/"
$ghDe   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fe2164f78855c6857b6d526e280e05b5f570e03/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/"

$id1 = "44e48f6d-14e0-46ea-9fd1-e0dacc693fab"
$id2 = "9108f6ff-b6e6-4f65-9bec-cc42006e03af"

$hash1 = "3de409930bc3354428a9bd2a0523794cf539c963"
$hash2 = "f8b59048bc8c7459296b55851d3372f6f026be07"

$status   = "Ready for handoff"
$overviewDate = "2016-37-18 02:37:22"

$zhDate1 = "2016-03-18 02:37:14"
$zhDate2 = "2016-03-18 02:37:14"
$deDate1 = "2016-03-18 02:37:22"
$deDate2 = "2016-03-18 02:37:22"

$md1 = "$id1.md"
$md2 = "$id2.md"

$xlf1zh = "$id1.$hash1.zh-cn.xlf"
$xlf2zh = "$id2.$hash2.zh-cn.xlf"
$xlf1de = "$id1.$hash1.de-de.xlf"
$xlf2de = "$id2.$hash2.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview" -> rows 4 and 5
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), ($ghSrc + $md1), "", "", $md1)
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), ($ghSrc + $md2), "", "", $md2)
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $overviewDate

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> rows 4 and 5
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($ghSrc + $md1), "", "", $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), ($ghSrc + $md1), "", "", ".md")
$wsZh.Range("C4").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), ($ghZh + $xlf1zh), "", "", $xlf1zh)
$wsZh.Range("E4").Value = $zhDate1
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($ghSrc + $md2), "", "", $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), ($ghSrc + $md2), "", "", ".md")
$wsZh.Range("C5").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), ($ghZh + $xlf2zh), "", "", $xlf2zh)
$wsZh.Range("E5").Value = $zhDate2
$wsZh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de" -> rows 4 and 5
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($ghSrc + $md1), "", "", $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), ($ghSrc + $md1), "", "", ".md")
$wsDe.Range("C4").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), ($ghDe + $xlf1de), "", "", $xlf1de)
$wsDe.Range("E4").Value = $deDate1
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($ghSrc + $md2), "", "", $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), ($ghSrc + $md2), "", "", ".md")
$wsDe.Range("C5").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), ($ghDe + $xlf2de), "", "", $xlf2de)
$wsDe.Range("E5").Value = $deDate2
$wsDe.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"
